# Replace the customer records in the "Direct Deposit" sheet with the
# restructured data set (classification validation / future-lab dataset).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("C001", "Lauren",   "Smith",     212524838, 35928397, "lauren.smith@contoso.com"),
    @("C002", "Robert",   "Rodriguez", 310308176, 10080973, "robert.rodriguez@contoso.com"),
    @("C003", "Lauren",   "Smith",     661028505, 69775754, "lauren.smith@gmail.com"),
    @("C004", "Amanda",   "Wilson",    708747422, 98891132, "amanda.wilson@contoso.com"),
    @("C005", "Jessica",  "Moore",     344512868, 23012863, "jessica.moore@contoso.com"),
    @("C006", "Lauren",   "Davis",     42683438,  85398874, "lauren.davis@hotmail.com"),
    @("C007", "Emily",    "Davis",     36667545,  64007205, "emily.davis@contoso.com"),
    @("C008", "Jennifer", "Gonzalez",  770642858, 94141095, "jennifer.gonzalez@contoso.com"),
    @("C009", "Michael",  "Taylor",    843461418, 89235880, "michael.taylor@icloud.com"),
    @("C010", "Michelle", "Taylor",    256177166, 54917528, "michelle.taylor@contoso.com"),
    @("C011", "Sarah",    "Wilson",    152723520, 25703299, "sarah.wilson@contoso.com"),
    @("C012", "James",    "Garcia",    116057488, 54053051, "james.garcia@hotmail.com")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $record = $data[$i]
    $ws.Cells.Item($row, 1).Value = $record[0]
    $ws.Cells.Item($row, 2).Value = $record[1]
    $ws.Cells.Item($row, 3).Value = $record[2]
    $ws.Cells.Item($row, 4).Value = $record[3]
    $ws.Cells.Item($row, 5).Value = $record[4]
    $ws.Cells.Item($row, 6).Value = $record[5]
}

# Column F ("Email") is a bit narrower now that the data set has changed.
# (29.15 round-trips to the stored OOXML width of exactly 30.)
$ws.Columns.Item(6).ColumnWidth = 29.15
